# Add the "mistral" sheet after the existing sheet and populate it with the
# same layout as the other per-model timing sheet: column A holds the run
# labels (bold, bordered, top/center aligned), column B holds the header in
# B1 plus the numeric timings in B2:B24.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "mistral"

$ws.Cells.Item(1, 2).Value = 'mistral'
$ws.Cells.Item(2, 1).Value = 'preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse'
$ws.Cells.Item(2, 2).Value = 3.820181445594042
$ws.Cells.Item(3, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_lenNone_gblFalse'
$ws.Cells.Item(3, 2).Value = 0.9920371054360548
$ws.Cells.Item(4, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_lenNone_gblFalse'
$ws.Cells.Item(4, 2).Value = 1.451674235796612
$ws.Cells.Item(5, 1).Value = 'preds_ns10_ws32_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse'
$ws.Cells.Item(5, 2).Value = 4.511714891931707
$ws.Cells.Item(6, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_qcache_lenNone_gblFalse'
$ws.Cells.Item(6, 2).Value = 0.8876167544943719
$ws.Cells.Item(7, 1).Value = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse'
$ws.Cells.Item(7, 2).Value = 4.646624851453602
$ws.Cells.Item(8, 1).Value = 'preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_sum_fused_rerun_lenNone_gblFalse'
$ws.Cells.Item(8, 2).Value = 8.122498173586132
$ws.Cells.Item(9, 1).Value = 'preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse'
$ws.Cells.Item(9, 2).Value = 8.10217903827626
$ws.Cells.Item(10, 1).Value = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_qcache_lenNone_gblFalse'
$ws.Cells.Item(10, 2).Value = 0.7907740784882081
$ws.Cells.Item(11, 1).Value = 'preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse'
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(12, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_snapkv_opt_lenNone_gblFalse'
$ws.Cells.Item(12, 2).Value = 0.654116009700185
$ws.Cells.Item(13, 1).Value = 'preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse'
$ws.Cells.Item(13, 2).Value = 3.176511576234317
$ws.Cells.Item(14, 1).Value = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse'
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(15, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_opt_qcache_lenNone_gblFalse'
$ws.Cells.Item(15, 2).Value = 1.695335097976731
$ws.Cells.Item(16, 1).Value = 'preds_ns5_ws200_mc1000_ea1.0_snks0_hopf_False_type_max_fused_lenNone_gblFalse'
$ws.Cells.Item(16, 2).Value = 0.6459789344721453
$ws.Cells.Item(17, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse'
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(18, 1).Value = 'preds_ns5_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse'
$ws.Cells.Item(18, 2).Value = 1.57627808441923
$ws.Cells.Item(19, 1).Value = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_snapkv_prof_lenNone_gblFalse'
$ws.Cells.Item(19, 2).Value = 1.416201545849141
$ws.Cells.Item(20, 1).Value = 'preds_ns5_ws200_mc2000_ea1.0_snks0_hopf_True_type_max_fused_opt_qcache_lenNone_gblFalse'
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(21, 1).Value = 'preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_snapkv_rerun_lenNone_gblFalse'
$ws.Cells.Item(21, 2).Value = 2.225652065976954
$ws.Cells.Item(22, 1).Value = 'preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_prof_lenNone_gblFalse'
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(23, 1).Value = 'preds_ns1_ws200_mc750_ea1.0_snks0_hopf_True_type_max_fused_prof_qcache_lenNone_gblFalse'
$ws.Cells.Item(23, 2).Value = 1.748998216443056
$ws.Cells.Item(24, 1).Value = 'preds_ns1_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_prof_qcache_lenNone_gblFalse'
$ws.Cells.Item(24, 2).Value = 0.1295784779876364


# Match the header/label styling used on the "llama3.1-8b-instruct" sheet:
# bold text, thin box border, centered horizontally and top-aligned vertically.
$labelRange = $ws.Range("A2:A24")
$labelRange.Font.Bold = $true
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4160
$labelRange.Borders.LineStyle = 1

$headerCell = $ws.Range("B1")
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108
$headerCell.VerticalAlignment = -4160
$headerCell.Borders.LineStyle = 1
